$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark Visualizations (row 3), lab01 (row 21), and lab02 (row 22) as rendered ("link_it" = TRUE)
$ws.Range("C3").Value = $true
$ws.Range("C21").Value = $true
$ws.Range("C22").Value = $true

# lab02 now has a topic: "R basics"
$ws.Range("D22").Value = "R basics"

# Update selection to match the final cursor position recorded in the workbook
$ws.Range("C23").Select()
